$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("展览").Range("F3").Value = 632
$wb.Worksheets.Item("展览").Range("F6").Value = 2819
$wb.Worksheets.Item("展览").Range("F8").Value = 55
$wb.Worksheets.Item("展览").Range("F10").Value = 580
$wb.Worksheets.Item("展览").Range("F11").Value = 22
$wb.Worksheets.Item("展览").Range("F12").Value = 323
$wb.Worksheets.Item("展览").Range("F14").Value = 5958
$wb.Worksheets.Item("展览").Range("F16").Value = 1047
$wb.Worksheets.Item("展览").Range("F17").Value = 8
$wb.Worksheets.Item("展览").Range("F18").Value = 231
$wb.Worksheets.Item("展览").Range("F19").Value = 171
$wb.Worksheets.Item("展览").Range("F20").Value = 88
$wb.Worksheets.Item("展览").Range("F21").Value = 534
$wb.Worksheets.Item("展览").Range("F22").Value = 29
$wb.Worksheets.Item("展览").Range("F23").Value = 29
$wb.Worksheets.Item("展览").Range("F24").Value = 95
$wb.Worksheets.Item("展览").Range("F25").Value = 1316
$wb.Worksheets.Item("展览").Range("F27").Value = 7
$wb.Worksheets.Item("展览").Range("F28").Value = 38
$wb.Worksheets.Item("展览").Range("F29").Value = 2058
$wb.Worksheets.Item("展览").Range("F30").Value = 180
$wb.Worksheets.Item("展览").Range("F31").Value = 351
$wb.Worksheets.Item("展览").Range("F33").Value = 3294
$wb.Worksheets.Item("演出").Range("F7").Value = 346
$wb.Worksheets.Item("演出").Range("F8").Value = 82
$wb.Worksheets.Item("演出").Range("F15").Value = 1005
$wb.Worksheets.Item("演出").Range("F18").Value = 629
$wb.Worksheets.Item("演出").Range("F22").Value = 350
$wb.Worksheets.Item("演出").Range("F23").Value = 301
$wb.Worksheets.Item("演出").Range("F24").Value = 4053
$wb.Worksheets.Item("演出").Range("F28").Value = 141
$wb.Worksheets.Item("演出").Range("F37").Value = 12
$wb.Worksheets.Item("演出").Range("G13").Value = 0
$wb.Worksheets.Item("本地生活").Range("F2").Value = 1802
$wb.Worksheets.Item("本地生活").Range("F5").Value = 2597
$wb.Worksheets.Item("本地生活").Range("F8").Value = 1492
$wb.Worksheets.Item("本地生活").Range("F12").Value = 641
$wb.Worksheets.Item("全部类型").Range("F2").Value = 1802
$wb.Worksheets.Item("全部类型").Range("F4").Value = 2597
$wb.Worksheets.Item("全部类型").Range("F7").Value = 1492
$wb.Worksheets.Item("全部类型").Range("F11").Value = 632
$wb.Worksheets.Item("全部类型").Range("F12").Value = 2819
$wb.Worksheets.Item("全部类型").Range("F13").Value = 55
$wb.Worksheets.Item("全部类型").Range("F14").Value = 641
$wb.Worksheets.Item("全部类型").Range("F15").Value = 580
$wb.Worksheets.Item("全部类型").Range("F16").Value = 82
$wb.Worksheets.Item("全部类型").Range("F17").Value = 22
$wb.Worksheets.Item("全部类型").Range("F18").Value = 323
$wb.Worksheets.Item("全部类型").Range("F20").Value = 5958
$wb.Worksheets.Item("全部类型").Range("F23").Value = 1047
$wb.Worksheets.Item("全部类型").Range("F24").Value = 231
$wb.Worksheets.Item("全部类型").Range("F25").Value = 171
$wb.Worksheets.Item("全部类型").Range("F26").Value = 88
$wb.Worksheets.Item("全部类型").Range("F27").Value = 534
$wb.Worksheets.Item("全部类型").Range("F32").Value = 29
$wb.Worksheets.Item("全部类型").Range("F34").Value = 350
$wb.Worksheets.Item("全部类型").Range("F35").Value = 301
$wb.Worksheets.Item("全部类型").Range("F39").Value = 141
$wb.Worksheets.Item("全部类型").Range("F41").Value = 38
$wb.Worksheets.Item("全部类型").Range("F44").Value = 2058
$wb.Worksheets.Item("全部类型").Range("F47").Value = 180
$wb.Worksheets.Item("全部类型").Range("F48").Value = 351
$wb.Worksheets.Item("全部类型").Range("F50").Value = 3294
